$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Replace-Text "2025-11-19 Wednesday" "2025-11-20 Thursday"
Replace-Text "505×2=1010" "723×7=5061"
Replace-Text "134×8=1072" "264×5=1320"
Replace-Text "571×9=5139" "719×9=6471"
Replace-Text "475×3=1425" "170×2=340"
Replace-Text "434×9=3906" "905×2=1810"
Replace-Text "573×8=4584" "512×8=4096"
Replace-Text "460×5=2300" "447×3=1341"
Replace-Text "368×5=1840" "901×4=3604"
Replace-Text "456×4=1824" "367×6=2202"
Replace-Text "293×9=2637" "159×7=1113"
Replace-Text "718×2=1436" "712×4=2848"
Replace-Text "293×6=1758" "501×9=4509"
Replace-Text "864×3=2592" "597×6=3582"
Replace-Text "380×8=3040" "286×2=572"
Replace-Text "471×6=2826" "197×6=1182"
Replace-Text "955×9=8595" "459×5=2295"
Replace-Text "839×9=7551" "869×6=5214"
Replace-Text "566×6=3396" "383×7=2681"
Replace-Text "625×3=1875" "252×4=1008"
Replace-Text "620×4=2480" "193×8=1544"
Replace-Text "872×8=6976" "482×9=4338"
Replace-Text "623×6=3738" "867×7=6069"
Replace-Text "896×4=3584" "647×7=4529"
Replace-Text "892×7=6244" "657×7=4599"
Replace-Text "558×4=2232" "313×2=626"

Write-Host "Done"
